# Updated cryptos list on Sat Jul 15 17:59:39 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# active worksheet for rows 2-51 to the latest scraped values.
#
# Note: several of the new Price values are plain decimal numbers (e.g.
# "1.001", "250.84"). Excel's Range.Value setter auto-detects such strings
# as numeric literals, which would silently reformat them (and introduce
# floating-point noise) instead of keeping them as literal text like the
# rest of the sheet. To keep these as text we prefix them with a leading
# apostrophe, exactly like typing `'1.001` into a cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.325.20"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "1.933.64"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'250.84"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'0.7141"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'0.3294"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").Value = "'27.56"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.07236"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").Value = "'0.8044"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").Value = "'0.08099"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "1.929.40"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").Value = "'5.464"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'94.49"
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "30.307.35"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "'252.75"
$ws.Range("E18").Value = "  -5.53%  "
$ws.Range("D19").Value = "'0.000008168"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'5.786"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("D21").Value = "2.184.95"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'6.968"
$ws.Range("D25").Value = "'9.735"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "'165.59"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").Value = "'2.343"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'19.28"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "'0.1287"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("D30").Value = "'1.352"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'4.416"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "'4.169"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").Value = "'0.05190"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").Value = "'1.259"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "'0.7448"
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.01962"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "'2.809"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "'78.62"
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("D41").Value = "'6.420"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "'0.4516"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "'2.015"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("D44").Value = "'0.8441"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").Value = "'9.721"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("D48").Value = "'7.426"
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("D49").Value = "'36.64"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").Value = "'0.4158"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").Value = "'0.06032"
$ws.Range("E51").Value = "  +0.20%  "

